$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price (D) and Volume(1h) (E) columns for rows 2-51
# so that numeric-looking strings (e.g. "1.033", "0.07389") are preserved exactly
# as text instead of being coerced into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.544.26"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "1.849.10"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "1.033"
$ws.Range("E4").Value = "  +2.59%  "
$ws.Range("D5").Value = "321.16"
$ws.Range("E5").Value = "  +3.56%  "
$ws.Range("D6").Value = "1.029"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").Value = "0.4373"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "0.3764"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "0.07389"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "0.8727"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "21.39"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "1.862.87"
$ws.Range("E12").Value = "  -7.77%  "
$ws.Range("D13").Value = "5.507"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "6.673"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "0.07212"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").Value = "82.55"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "1.035"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "0.000009034"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "1.029"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "15.38"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "27.548.10"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "5.237"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").Value = "11.34"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").Value = "2.071.74"
$ws.Range("E24").Value = "  -7.98%  "
$ws.Range("D25").Value = "157.49"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").Value = "1.922"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").Value = "5.253"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "1.956"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("D30").Value = "116.65"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "0.09024"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "0.7597"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").Value = "1.192"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "4.492"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "2.878"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").Value = "1.030"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").Value = "1.149"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("D38").Value = "0.01970"
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("D39").Value = "0.05279"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "0.5141"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "2.798"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").Value = "0.1670"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "6.689"
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("D44").Value = "8.457"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "108.84"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").Value = "10.49"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "1.704"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "0.06402"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "0.4633"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "1.853"
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("D51").Value = "39.12"
$ws.Range("E51").Value = "  +4.12%  "
